# Generate Report for Handoff
# Swap the two tracked files' row order/content on all three sheets
# ("84185a1c..." now in sync, "0b3b24c3..." now ready for handoff with a stale-handback error),
# and widen the "Error Detail" column on the locale sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ccf5ab36660402f7de66c7fbe2295103885f8c2/e2e/0b3b24c3-c269-480f-93d9-20066a7100db.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ed95c1d05460115606b8a11754d4ed3d0b30e8e0/e2e/0b3b24c3-c269-480f-93d9-20066a7100db.md."

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.md"
$ws1.Range("B2").Value = "e2e\84185a1c-1e73-4263-9277-ac83868f3c3e.md"

$ws1.Range("A3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.md"
$ws1.Range("B3").Value = "e2e\0b3b24c3-c269-480f-93d9-20066a7100db.md"
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-17 00:44:47"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.md"
$ws2.Range("G2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.b61612364fae30cdf40972953b77135a4db84469.zh-cn.xlf"
$ws2.Range("I2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.md"
$ws2.Range("J2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.b61612364fae30cdf40972953b77135a4db84469.zh-cn.xlf"

$ws2.Range("A3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.8f482a754cd4f7a5915a93825b5e29bb8cb88cb6.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-17 00:44:42"
$ws2.Range("I3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.md"
$ws2.Range("J3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.8f482a754cd4f7a5915a93825b5e29bb8cb88cb6.zh-cn.xlf"
$ws2.Range("P3").Value = $errorDetail

$ws2.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.md"
$ws3.Range("G2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.b61612364fae30cdf40972953b77135a4db84469.de-de.xlf"
$ws3.Range("I2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.md"
$ws3.Range("J2").Value = "84185a1c-1e73-4263-9277-ac83868f3c3e.b61612364fae30cdf40972953b77135a4db84469.de-de.xlf"

$ws3.Range("A3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.8f482a754cd4f7a5915a93825b5e29bb8cb88cb6.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-17 00:44:47"
$ws3.Range("I3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.md"
$ws3.Range("J3").Value = "0b3b24c3-c269-480f-93d9-20066a7100db.8f482a754cd4f7a5915a93825b5e29bb8cb88cb6.de-de.xlf"
$ws3.Range("P3").Value = $errorDetail

$ws3.Columns.Item(16).ColumnWidth = 39.17
